$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue 'D2' '27.980.69'
Set-TextValue 'E2' '  +3.39%  '

Set-TextValue 'D3' '1.726.61'
Set-TextValue 'E3' '  +3.04%  '

Set-TextValue 'E4' '  -0.21%  '

Set-TextValue 'D5' '218.59'
Set-TextValue 'E5' '  +1.58%  '

Set-TextValue 'E6' '  +1.54%  '

Set-TextValue 'D8' '24.09'
Set-TextValue 'E8' '  +13.51%  '

Set-TextValue 'E9' '  +3.54%  '

Set-TextValue 'D10' '0.0635'
Set-TextValue 'E10' '  +2.20%  '

Set-TextValue 'E11' '  +2.14%  '

Set-TextValue 'D12' '1.970.49'
Set-TextValue 'E12' '  +3.05%  '

Set-TextValue 'D13' '1.723.33'
Set-TextValue 'E13' '  +2.46%  '

Set-TextValue 'E14' '  +3.71%  '

Set-TextValue 'D15' '0.568'
Set-TextValue 'E15' '  +6.19%  '

Set-TextValue 'D16' '67.92'
Set-TextValue 'E16' '  +2.90%  '

Set-TextValue 'D17' '27.926.17'
Set-TextValue 'E17' '  +3.24%  '

Set-TextValue 'D18' '244.05'
Set-TextValue 'E18' '  +3.02%  '

Set-TextValue 'D19' '0.0₃0757'
Set-TextValue 'E19' '  +2.32%  '

Set-TextValue 'D20' '7.87'
Set-TextValue 'E20' '  -3.32%  '

Set-TextValue 'D22' '4.66'
Set-TextValue 'E22' '  +4.39%  '

Set-TextValue 'D23' '9.77'
Set-TextValue 'E23' '  +4.85%  '

Set-TextValue 'E24' '  +0.84%  '

Set-TextValue 'D25' '149.50'
Set-TextValue 'E25' '  +2.37%  '

Set-TextValue 'D26' '7.53'
Set-TextValue 'E26' '  +4.27%  '

Set-TextValue 'D27' '16.85'
Set-TextValue 'E27' '  +3.10%  '

Set-TextValue 'E28' '  +1.84%  '

Set-TextValue 'D29' '1.00'
Set-TextValue 'E29' '  -0.33%  '

Set-TextValue 'E30' '  +2.88%  '

Set-TextValue 'E31' '  +1.91%  '

Set-TextValue 'E32' '  +2.80%  '

Set-TextValue 'E33' '  +3.76%  '

Set-TextValue 'D34' '1.489.52'
Set-TextValue 'E34' '  -3.06%  '

Set-TextValue 'E35' '  -2.09%  '

Set-TextValue 'B36' 'ImmutableX'
Set-TextValue 'C36' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D36' '0.613'
Set-TextValue 'E36' '  +3.13%  '

Set-TextValue 'B37' 'ARBITRUM'
Set-TextValue 'C37' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D37' '0.960'
Set-TextValue 'E37' '  +4.76%  '

Set-TextValue 'E38' '  +0.60%  '

Set-TextValue 'E39' '  +0.73%  '

Set-TextValue 'E40' '  +0.33%  '

Set-TextValue 'D41' '71.63'
Set-TextValue 'E41' '  +6.01%  '

Set-TextValue 'D42' '5.85'
Set-TextValue 'E42' '  +5.51%  '

Set-TextValue 'B44' 'RocketPoolETH'
Set-TextValue 'C44' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 'D44' '1.874.77'
Set-TextValue 'E44' '  +3.10%  '

Set-TextValue 'B45' 'MXToken'
Set-TextValue 'C45' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D45' '2.29'
Set-TextValue 'E45' '  +1.42%  '

Set-TextValue 'D46' '0.795'
Set-TextValue 'E46' '  +1.61%  '

Set-TextValue 'D47' '1.75'
Set-TextValue 'E47' '  +12.32%  '

Set-TextValue 'D48' '91.36'
Set-TextValue 'E48' '  +0.56%  '

Set-TextValue 'E49' '  +3.45%  '

Set-TextValue 'D50' '0.105'
Set-TextValue 'E50' '  +1.28%  '

Set-TextValue 'D51' '8.24'
Set-TextValue 'E51' '  +2.56%  '

